$d = $word.ActiveDocument

# Update the title: "Multiobjective" -> "Parking", "Optimization" -> "Garage Case Study"
$d.Content.Find.Execute("Multiobjective", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Parking", 2)

$d.Content.Find.Execute("Optimization", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Garage Case Study", 2)

# Update the date: "Mar. 1" -> "Mar. 15"
$d.Content.Find.Execute("Mar. 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mar. 15", 2)
